$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 720345.1
$ws.Range("I28").Value = 2501283
$ws.Range("J28").Value = 7970
$ws.Range("K28").Value = 2501283
$ws.Range("L28").Value = 7970
$ws.Range("M28").Value = -2500798
$ws.Range("N28").Value = -8940
$ws.Range("H41").Value = 269.25
$ws.Range("I41").Value = 236
$ws.Range("K41").Value = 236
$ws.Range("M41").Value = 204
$ws.Range("H70").Value = 2021.75
$ws.Range("J70").Value = 2021.75
$ws.Range("L70").Value = 6065.25
$ws.Range("N70").Value = -6605.25
$ws.Range("H73").Value = 2021.75
$ws.Range("J73").Value = 2021.75
$ws.Range("L73").Value = 6065.25
$ws.Range("N73").Value = -7937.25
$ws.Range("H76").Value = 838822.7
$ws.Range("I76").Value = 1253898.8
$ws.Range("K76").Value = 1253898.8
$ws.Range("M76").Value = -1253583.8
$ws.Range("H79").Value = 838822.7
$ws.Range("I79").Value = 1253898.8
$ws.Range("K79").Value = 1253898.8
$ws.Range("M79").Value = -1252806.8
$ws.Range("H86").Value = 1004325.1
$ws.Range("I86").Value = 1114805.2
$ws.Range("K86").Value = 1114805.2
$ws.Range("M86").Value = -1113682.2
$ws.Range("H88").Value = 1233.4445
$ws.Range("J88").Value = 1435.8572
$ws.Range("L88").Value = 1435.8572
$ws.Range("N88").Value = -2247.8572
$ws.Range("H89").Value = 1004325.1
$ws.Range("I89").Value = 1114805.2
$ws.Range("K89").Value = 5574026
$ws.Range("M89").Value = -5568410
$ws.Range("H91").Value = 1233.4445
$ws.Range("J91").Value = 1435.8572
$ws.Range("L91").Value = 1435.8572
$ws.Range("N91").Value = -4243.8572
$ws.Range("H98").Value = 4159.381
$ws.Range("I98").Value = 4365.5386
$ws.Range("J98").Value = 3824.375
$ws.Range("K98").Value = 4365.5386
$ws.Range("L98").Value = 3824.375
$ws.Range("M98").Value = -2867.5386
$ws.Range("N98").Value = -6820.375
$ws.Range("H100").Value = 4626
$ws.Range("I100").Value = 3078.3572
$ws.Range("J100").Value = 5610.864
$ws.Range("K100").Value = 3078.3572
$ws.Range("L100").Value = 5610.864
$ws.Range("M100").Value = -2537.3572
$ws.Range("N100").Value = -6692.864
$ws.Range("H101").Value = 2985.8
$ws.Range("I101").Value = 3483.75
$ws.Range("K101").Value = 10451.25
$ws.Range("M101").Value = -8829.25
$ws.Range("H113").Value = 2884.2
$ws.Range("I113").Value = 2333.2222
$ws.Range("J113").Value = 3335
$ws.Range("K113").Value = 2333.2222
$ws.Range("L113").Value = 3335
$ws.Range("M113").Value = 920.7777999999998
$ws.Range("N113").Value = -9843
$ws.Range("H122").Value = 4159.381
$ws.Range("I122").Value = 4365.5386
$ws.Range("J122").Value = 3824.375
$ws.Range("K122").Value = 13096.6158
$ws.Range("L122").Value = 11473.125
$ws.Range("M122").Value = -10646.6158
$ws.Range("N122").Value = -16373.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6037.375
$ws.Range("I63").Value = 2449.5
$ws.Range("K63").Value = 2449.5
$ws.Range("M63").Value = -1763.5
$ws.Range("H66").Value = 6037.375
$ws.Range("I66").Value = 2449.5
$ws.Range("K66").Value = 12247.5
$ws.Range("M66").Value = -8815.5
$ws.Range("H97").Value = 166501060
$ws.Range("I97").Value = 249750600
$ws.Range("K97").Value = 249750600
$ws.Range("M97").Value = -249750104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 78276.30499999999
$ws.Range("I94").Value = 1199.75
$ws.Range("J94").Value = 201598.8
$ws.Range("K94").Value = 1199.75
$ws.Range("L94").Value = 201598.8
$ws.Range("M94").Value = -748.75
$ws.Range("N94").Value = -202500.8
$ws.Range("H107").Value = 16008.474
$ws.Range("I107").Value = 19867.785
$ws.Range("K107").Value = 19867.785
$ws.Range("M107").Value = -17947.785

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5444
$ws.Range("I99").Value = 5534.8335
$ws.Range("K99").Value = 5534.8335
$ws.Range("M99").Value = -4036.8335
$ws.Range("H126").Value = 5444
$ws.Range("I126").Value = 5534.8335
$ws.Range("K126").Value = 16604.5005
$ws.Range("M126").Value = -14134.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 63997.125
$ws.Range("J68").Value = 70017.38
$ws.Range("L68").Value = 210052.14
$ws.Range("N68").Value = -211674.14
$ws.Range("H71").Value = 63997.125
$ws.Range("J71").Value = 70017.38
$ws.Range("L71").Value = 630156.42
$ws.Range("N71").Value = -638268.42
$ws.Range("H118").Value = 9741.125
$ws.Range("I118").Value = 7305.8
$ws.Range("K118").Value = 21917.4
$ws.Range("M118").Value = -20674.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4397.4883
$ws.Range("I102").Value = 2588.8
$ws.Range("J102").Value = 5366.4287
$ws.Range("K102").Value = 2588.8
$ws.Range("L102").Value = 5366.4287
$ws.Range("M102").Value = -966.8000000000002
$ws.Range("N102").Value = -8610.4287
$ws.Range("H107").Value = 241.66667
$ws.Range("I107").Value = 241.66667
$ws.Range("K107").Value = 241.66667
$ws.Range("M107").Value = 1678.33333
$ws.Range("H122").Value = 3934.75
$ws.Range("I122").Value = 3870.76
$ws.Range("J122").Value = 4080.182
$ws.Range("K122").Value = 11612.28
$ws.Range("L122").Value = 12240.546
$ws.Range("M122").Value = -9162.280000000001
$ws.Range("N122").Value = -17140.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9950
$ws.Range("J40").Value = 9950
$ws.Range("L40").Value = 9950
$ws.Range("N40").Value = -10222
$ws.Range("H46").Value = 1307.4595
$ws.Range("I46").Value = 955.8889
$ws.Range("K46").Value = 955.8889
$ws.Range("M46").Value = -767.8889
$ws.Range("H55").Value = 2021.6666
$ws.Range("I55").Value = 524.7646999999999
$ws.Range("J55").Value = 5657
$ws.Range("K55").Value = 524.7646999999999
$ws.Range("L55").Value = 5657
$ws.Range("M55").Value = -351.7646999999999
$ws.Range("N55").Value = -6003
$ws.Range("H61").Value = 6456.8667
$ws.Range("I61").Value = 5192.5386
$ws.Range("J61").Value = 14675
$ws.Range("K61").Value = 5192.5386
$ws.Range("L61").Value = 14675
$ws.Range("M61").Value = -4990.5386
$ws.Range("N61").Value = -15079
$ws.Range("H113").Value = 6456.8667
$ws.Range("I113").Value = 5192.5386
$ws.Range("J113").Value = 14675
$ws.Range("K113").Value = 5192.5386
$ws.Range("L113").Value = 14675
$ws.Range("M113").Value = -3022.5386
$ws.Range("N113").Value = -19015
$ws.Range("H122").Value = 6101.857
$ws.Range("J122").Value = 7005.0625
$ws.Range("L122").Value = 21015.1875
$ws.Range("N122").Value = -25915.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 19197.4
$ws.Range("I62").Value = 13309.25
$ws.Range("K62").Value = 13309.25
$ws.Range("M62").Value = -12685.25
$ws.Range("H65").Value = 19197.4
$ws.Range("I65").Value = 13309.25
$ws.Range("K65").Value = 66546.25
$ws.Range("M65").Value = -63426.25
$ws.Range("H113").Value = 1335.3529
$ws.Range("I113").Value = 1335.3529
$ws.Range("K113").Value = 4006.0587
$ws.Range("M113").Value = -1836.0587
$ws.Range("H122").Value = 3743.5854
$ws.Range("I122").Value = 3508.7878
$ws.Range("K122").Value = 10526.3634
$ws.Range("M122").Value = -8076.3634
